$wb = $excel.ActiveWorkbook

# --- System Settings sheet: insert row 9 with NCAP_AF entry ---
$sys = $wb.Worksheets.Item("System Settings")
$sys.Rows("9:9").Insert() | Out-Null
$sys.Range("C9").Value = "NCAP_AF"
$sys.Range("D9").Value = 0
$sys.Range("E9").Value = 3
$sys.Range("C11").Select() | Out-Null

# --- fuels sheet: add new side table (M3:Q6), mirroring the ~FI_Comm table ---
$fuels = $wb.Worksheets.Item("fuels")
$fuels.Range("B3:F4").Copy($fuels.Range("M3:Q4")) | Out-Null

# fix header labels/order for the new table
$fuels.Range("O4").Value = "timeslicelevel"
$fuels.Range("P4").Value = "Unit"
$fuels.Range("Q4").Value = "description"

$fuels.Range("M5").Value = "NRG"
$fuels.Range("N5").Value = "ELC_Sol-IND"
$fuels.Range("O5").Value = "DAYNITE"
$fuels.Range("P5").Value = "TWh"
$fuels.Range("Q5").Value = "Solar electricity produced in - India"

$fuels.Range("M6").Value = "NRG"
$fuels.Range("N6").Value = "ELC_Win-IND"
$fuels.Range("O6").Value = "DAYNITE"
$fuels.Range("P6").Value = "TWh"
$fuels.Range("Q6").Value = "Wind electricity produced in - India"

$fuels.Select() | Out-Null
$fuels.Range("M5").Select() | Out-Null
